# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E26) cycled through period codes
# 1801, 1712, 1711, ..., 1703 (newest-first, with 1801 leading).
# The update rotates the list so the new period set 1703..1712 comes
# first (part 1 of the new account statements) followed by 1801, and
# the "Valor Mora" (F column) value travels together with its period
# label (F16 <-> F26 swap follows E16/E26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E), rows 16-26: new period ordering.
$ws.Range("E16").Value = "1703"
$ws.Range("E17").Value = "1704"
$ws.Range("E18").Value = "1705"
$ws.Range("E19").Value = "1706"
$ws.Range("E20").Value = "1707"
$ws.Range("E21").Value = "1708"
$ws.Range("E22").Value = "1709"
$ws.Range("E23").Value = "1710"
$ws.Range("E24").Value = "1711"
$ws.Range("E25").Value = "1712"
$ws.Range("E26").Value = "1801"

# Valor Mora (column F) follows its period label: 1703 now carries the
# value previously on 1801's row (29509) and vice versa.
$ws.Range("F16").Value = 29509
$ws.Range("F26").Value = 28526
